$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert the new credit row (row 5) for the explosion sprite sheet asset.
$ws.Range("A5").Value = "explosion_01_strip13.png"
$ws.Range("D5").Value = "`"Bleed - http://remusprites.carbonmade.com/`""
$ws.Range("C5").Value = "Creative Commons License 3"
$ws.Range("B5").Value = "https://opengameart.org/content/simple-explosion-bleeds-game-art"

# Style the new Notes/Other cell like the other note cells: Arial 14pt,
# but with the automatic (theme) text color rather than the custom gray.
$ws.Range("D5").Font.Name = "Arial"
$ws.Range("D5").Font.Size = 14
$ws.Range("D5").Font.ThemeColor = 1

$ws.Rows.Item(5).RowHeight = 18

# Update the saved selection/scroll state to match the new view.
$ws.Range("B8").Select()
